$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Fill in row 8 with a new data row (mirrors the pattern used in rows 2-7)
$ws.Range("A8").Value = "TabStats_V01"
$ws.Range("B8").Value = "F2"
$ws.Range("C8").Value = 10
$ws.Range("D8").Value = 30
$ws.Range("E8").Value = 2
# F8 needs to hold the numeric-looking value "635723.112" as text (not a
# number) without altering the cell's style, so compute it as text via a
# scratch formula cell and paste the resulting value back in.
$ws.Range("Z1").Formula = '=TEXT(635723.112,"0.000")'
$ws.Range("Z1").Copy()
$ws.Range("F8").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("Z1").ClearContents()
$ws.Range("G8").Value = 435

# Update the active selection to D8, matching the saved view state
$ws.Range("D8").Select()
